$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.49677774085708
$ws.Range("C2").Value = 8.006696187756392
$ws.Range("D2").Value = 6.067008707835074
$ws.Range("E2").Value = 8.626259165873787
$ws.Range("G2").Value = 3.761216909668637
$ws.Range("I2").Value = 37.00665040959698
$ws.Range("K2").Value = 17.73292023895231
$ws.Range("L2").Value = 10.82127092317198
$ws.Range("N2").Value = 24.88396329244865

$ws.Range("B3").Value = 19.32935816687735
$ws.Range("C3").Value = 7.797926514967203
$ws.Range("D3").Value = 5.962468258102314
$ws.Range("E3").Value = 8.639204754926084
$ws.Range("G3").Value = 3.764870437854596
$ws.Range("I3").Value = 36.93209736109324
$ws.Range("K3").Value = 17.60785762403586
$ws.Range("L3").Value = 10.81372361085618
$ws.Range("N3").Value = 24.89150680294505

$ws.Range("B4").Value = 19.23152864541431
$ws.Range("C4").Value = 7.669471581198842
$ws.Range("D4").Value = 5.899276222638383
$ws.Range("E4").Value = 8.648336106577288
$ws.Range("G4").Value = 3.767229672399893
$ws.Range("I4").Value = 36.89201518362711
$ws.Range("K4").Value = 17.53550469804664
$ws.Range("L4").Value = 10.81118369423664
$ws.Range("N4").Value = 24.89786352746364

$ws.Range("B5").Value = 19.19295131514418
$ws.Range("C5").Value = 7.617143538751435
$ws.Range("D5").Value = 5.873811676034959
$ws.Range("E5").Value = 8.652354780387183
$ws.Range("G5").Value = 3.768220347385175
$ws.Range("I5").Value = 36.87711855432887
$ws.Range("K5").Value = 17.50716184931671
$ws.Range("L5").Value = 10.81067609667224
$ws.Range("N5").Value = 24.90088736660482

$ws.Range("B6").Value = 19.18662454915604
$ws.Range("C6").Value = 7.608458195312032
$ws.Range("D6").Value = 5.869601676091388
$ws.Range("E6").Value = 8.653040055345333
$ws.Range("G6").Value = 3.76838661915145
$ws.Range("I6").Value = 36.87473192128057
$ws.Range("K6").Value = 17.50252520851809
$ws.Range("L6").Value = 10.81062368532181
$ws.Range("N6").Value = 24.90141563970197

$ws.Range("B7").Value = 19.23100310807154
$ws.Range("C7").Value = 7.668765676347943
$ws.Range("D7").Value = 5.898931590312362
$ws.Range("E7").Value = 8.648389098736363
$ws.Range("G7").Value = 3.767242914340293
$ws.Range("I7").Value = 36.89180845711666
$ws.Range("K7").Value = 17.53511780152862
$ws.Range("L7").Value = 10.81117471212938
$ws.Range("N7").Value = 24.89790255367391

$ws.Range("B8").Value = 19.43804544604204
$ws.Range("C8").Value = 7.934819017537318
$ws.Range("D8").Value = 6.030778198724686
$ws.Range("E8").Value = 8.630477444763834
$ws.Range("G8").Value = 3.762452646582763
$ws.Range("I8").Value = 36.97976636443872
$ws.Range("K8").Value = 17.68889423470445
$ws.Range("L8").Value = 10.81823455804441
$ws.Range("N8").Value = 24.8862060797712

$ws.Range("B9").Value = 19.8815823254304
$ws.Range("C9").Value = 8.450805186543295
$ws.Range("D9").Value = 6.295636461799411
$ws.Range("E9").Value = 8.604728986235276
$ws.Range("G9").Value = 3.753973952035885
$ws.Range("I9").Value = 37.19714974540545
$ws.Range("K9").Value = 18.02440417090795
$ws.Range("L9").Value = 10.84864678242047
$ws.Range("N9").Value = 24.87697218102321

$ws.Range("B10").Value = 20.22765704658688
$ws.Range("C10").Value = 8.821872199680032
$ws.Range("D10").Value = 6.491953033422772
$ws.Range("E10").Value = 8.591516139161065
$ws.Range("G10").Value = 3.748295422338721
$ws.Range("I10").Value = 37.38380297191902
$ws.Range("K10").Value = 18.28981277681936
$ws.Range("L10").Value = 10.88100955207833
$ws.Range("N10").Value = 24.87856466934975

$ws.Range("B11").Value = 20.38888091379192
$ws.Range("C11").Value = 8.988076075180759
$ws.Range("D11").Value = 6.58120152784191
$ws.Range("E11").Value = 8.586740946662562
$ws.Range("G11").Value = 3.745830201122763
$ws.Range("I11").Value = 37.47445370564081
$ws.Range("K11").Value = 18.41424691458913
$ws.Range("L11").Value = 10.89788410470316
$ws.Range("N11").Value = 24.88111230661305

$ws.Range("B12").Value = 20.4504254127393
$ws.Range("C12").Value = 9.050575152996828
$ws.Range("D12").Value = 6.614954599779884
$ws.Range("E12").Value = 8.585110033472935
$ws.Range("G12").Value = 3.744913534450253
$ws.Range("I12").Value = 37.50959492534092
$ws.Range("K12").Value = 18.46186220994365
$ws.Range("L12").Value = 10.90458103980262
$ws.Range("N12").Value = 24.88233940742102

$ws.Range("B13").Value = 20.43714964689871
$ws.Range("C13").Value = 9.037135354558892
$ws.Range("D13").Value = 6.607687747352678
$ws.Range("E13").Value = 8.585453397348582
$ws.Range("G13").Value = 3.745110206872133
$ws.Range("I13").Value = 37.50199064194185
$ws.Range("K13").Value = 18.45158599728607
$ws.Range("L13").Value = 10.90312513031455
$ws.Range("N13").Value = 24.88206345797117

$ws.Range("B14").Value = 20.39393465092336
$ws.Range("C14").Value = 8.993227059971522
$ws.Range("D14").Value = 6.583979479738182
$ws.Range("E14").Value = 8.586603217968785
$ws.Range("G14").Value = 3.745754449122114
$ws.Range("I14").Value = 37.47732855741443
$ws.Range("K14").Value = 18.41815455257785
$ws.Range("L14").Value = 10.89842893083848
$ws.Range("N14").Value = 24.88120800181837

$ws.Range("B15").Value = 20.36752677893179
$ws.Range("C15").Value = 8.966272978428821
$ws.Range("D15").Value = 6.56945080891988
$ws.Range("E15").Value = 8.5873306028004
$ws.Range("G15").Value = 3.746151258548199
$ws.Range("I15").Value = 37.46232792513214
$ws.Range("K15").Value = 18.39774014896477
$ws.Range("L15").Value = 10.89559225829011
$ws.Range("N15").Value = 24.88071818282195

$ws.Range("B16").Value = 20.21719250297452
$ws.Range("C16").Value = 8.810952429879732
$ws.Range("D16").Value = 6.486116247515648
$ws.Range("E16").Value = 8.591853043747731
$ws.Range("G16").Value = 3.74845889523217
$ws.Range("I16").Value = 37.37799327357541
$ws.Range("K16").Value = 18.28175198965133
$ws.Range("L16").Value = 10.87994985794332
$ws.Range("N16").Value = 24.87843488027868

$ws.Range("B17").Value = 20.12590131342788
$ws.Range("C17").Value = 8.714955103381218
$ws.Range("D17").Value = 6.434953130722381
$ws.Range("E17").Value = 8.594943627617109
$ws.Range("G17").Value = 3.749904695862694
$ws.Range("I17").Value = 37.32771907140199
$ws.Range("K17").Value = 18.2115182125168
$ws.Range("L17").Value = 10.87090330468707
$ws.Range("N17").Value = 24.87750126637849

$ws.Range("B18").Value = 20.0737532458324
$ws.Range("C18").Value = 8.659499167034905
$ws.Range("D18").Value = 6.405522128710215
$ws.Range("E18").Value = 8.596837551446725
$ws.Range("G18").Value = 3.750747392326582
$ws.Range("I18").Value = 37.29934384639389
$ws.Range("K18").Value = 18.17147199892682
$ws.Range("L18").Value = 10.86590275347285
$ws.Range("N18").Value = 24.87713587691416

$ws.Range("B19").Value = 20.05616027220184
$ws.Range("C19").Value = 8.640683459438845
$ws.Range("D19").Value = 6.395557820703147
$ws.Range("E19").Value = 8.597498785000088
$ws.Range("G19").Value = 3.751034626256422
$ws.Range("I19").Value = 37.28982977240251
$ws.Range("K19").Value = 18.15797433352736
$ws.Range("L19").Value = 10.86424455657739
$ws.Range("N19").Value = 24.87704162821113

$ws.Range("B20").Value = 20.1355825170615
$ws.Range("C20").Value = 8.725199596404922
$ws.Range("D20").Value = 6.440400136025606
$ws.Range("E20").Value = 8.594602594962373
$ws.Range("G20").Value = 3.749749638821477
$ws.Range("I20").Value = 37.33301492015703
$ws.Range("K20").Value = 18.21895871873445
$ws.Range("L20").Value = 10.87184535612329
$ws.Range("N20").Value = 24.87758288977954

$ws.Range("B21").Value = 20.40661500029746
$ws.Range("C21").Value = 9.006136391661014
$ws.Range("D21").Value = 6.590944623963119
$ws.Range("E21").Value = 8.586260677162086
$ws.Range("G21").Value = 3.745564762839339
$ws.Range("I21").Value = 37.4845504236335
$ws.Range("K21").Value = 18.42796104939918
$ws.Range("L21").Value = 10.89980001153883
$ws.Range("N21").Value = 24.88145214820886

$ws.Range("B22").Value = 20.58659630246345
$ws.Range("C22").Value = 9.187161064215688
$ws.Range("D22").Value = 6.68906591174446
$ws.Range("E22").Value = 8.581842310140171
$ws.Range("G22").Value = 3.742927923826829
$ws.Range("I22").Value = 37.58832463163818
$ws.Range("K22").Value = 18.56742241186448
$ws.Range("L22").Value = 10.91985745495799
$ws.Range("N22").Value = 24.88551025523592

$ws.Range("B23").Value = 20.49029409928982
$ws.Range("C23").Value = 9.09080141819976
$ws.Range("D23").Value = 6.636732535155867
$ws.Range("E23").Value = 8.584106013557777
$ws.Range("G23").Value = 3.744326301900405
$ws.Range("I23").Value = 37.53250911313006
$ws.Range("K23").Value = 18.4927394462653
$ws.Range("L23").Value = 10.90898981502253
$ws.Range("N23").Value = 24.88320438353225

$ws.Range("B24").Value = 20.13120459286178
$ws.Range("C24").Value = 8.720568884000809
$ws.Range("D24").Value = 6.437937594108976
$ws.Range("E24").Value = 8.59475641104479
$ws.Range("G24").Value = 3.74981970432209
$ws.Range("I24").Value = 37.3306190210435
$ws.Range("K24").Value = 18.21559382854515
$ws.Range("L24").Value = 10.87141883015907
$ws.Range("N24").Value = 24.87754545403448

$ws.Range("B25").Value = 19.75784468698123
$ws.Range("C25").Value = 8.312317742654482
$ws.Range("D25").Value = 6.223531143743487
$ws.Range("E25").Value = 8.610691890689209
$ws.Range("G25").Value = 3.756170435846352
$ws.Range("I25").Value = 37.13357521376244
$ws.Range("K25").Value = 17.93018308473962
$ws.Range("L25").Value = 10.83865291799414
$ws.Range("N25").Value = 24.87800065292863
